$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2859515.5
$ws.Range("I62").Value = 3975360.8
$ws.Range("J62").Value = 19182.363
$ws.Range("K62").Value = 3975360.8
$ws.Range("L62").Value = 19182.363
$ws.Range("M62").Value = -3974736.8
$ws.Range("N62").Value = -20430.363
$ws.Range("H65").Value = 2859515.5
$ws.Range("I65").Value = 3975360.8
$ws.Range("J65").Value = 19182.363
$ws.Range("K65").Value = 19876804
$ws.Range("L65").Value = 95911.815
$ws.Range("M65").Value = -19873684
$ws.Range("N65").Value = -102151.815
$ws.Range("H92").Value = 695886.9399999999
$ws.Range("I92").Value = 1011419.44
$ws.Range("J92").Value = 1715.4
$ws.Range("K92").Value = 1011419.44
$ws.Range("L92").Value = 1715.4
$ws.Range("M92").Value = -1010171.44
$ws.Range("N92").Value = -4211.4
$ws.Range("H116").Value = 2311.1765
$ws.Range("I116").Value = 2449.2856
$ws.Range("J116").Value = 1666.6666
$ws.Range("K116").Value = 2449.2856
$ws.Range("L116").Value = 1666.6666
$ws.Range("M116").Value = 992.7143999999998
$ws.Range("N116").Value = -8550.6666
$ws.Range("H132").Value = 52056.15
$ws.Range("I132").Value = 52056.15
$ws.Range("K132").Value = 156168.45
$ws.Range("M132").Value = -153638.45
$ws.Range("H137").Value = 31251462
$ws.Range("I137").Value = 43479400
$ws.Range("J137").Value = 2288.7778
$ws.Range("K137").Value = 130438200
$ws.Range("L137").Value = 6866.3334
$ws.Range("M137").Value = -130435650
$ws.Range("N137").Value = -11966.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 698.5769
$ws.Range("I80").Value = 199.5
$ws.Range("J80").Value = 920.3889
$ws.Range("K80").Value = 199.5
$ws.Range("L80").Value = 920.3889
$ws.Range("M80").Value = 798.5
$ws.Range("N80").Value = -2916.3889
$ws.Range("H83").Value = 698.5769
$ws.Range("I83").Value = 199.5
$ws.Range("J83").Value = 920.3889
$ws.Range("K83").Value = 997.5
$ws.Range("L83").Value = 4601.944500000001
$ws.Range("M83").Value = 3994.5
$ws.Range("N83").Value = -14585.9445

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5443.653
$ws.Range("I31").Value = 1533.4642
$ws.Range("J31").Value = 10657.238
$ws.Range("K31").Value = 1533.4642
$ws.Range("L31").Value = 10657.238
$ws.Range("M31").Value = -1238.4642
$ws.Range("N31").Value = -11247.238
$ws.Range("H34").Value = 5443.653
$ws.Range("I34").Value = 1533.4642
$ws.Range("J34").Value = 10657.238
$ws.Range("K34").Value = 1533.4642
$ws.Range("L34").Value = 10657.238
$ws.Range("M34").Value = -1331.4642
$ws.Range("N34").Value = -11061.238
$ws.Range("H43").Value = 26625
$ws.Range("J43").Value = 26625
$ws.Range("L43").Value = 26625
$ws.Range("N43").Value = -26993
$ws.Range("H58").Value = 2082.5642
$ws.Range("I58").Value = 1153.3077
$ws.Range("J58").Value = 3941.077
$ws.Range("K58").Value = 1153.3077
$ws.Range("L58").Value = 3941.077
$ws.Range("M58").Value = -950.3077000000001
$ws.Range("N58").Value = -4347.077
$ws.Range("H92").Value = 29997.5
$ws.Range("J92").Value = 29997.5
$ws.Range("L92").Value = 29997.5
$ws.Range("N92").Value = -34989.5
$ws.Range("H95").Value = 24646.285
$ws.Range("J95").Value = 24646.285
$ws.Range("L95").Value = 24646.285
$ws.Range("N95").Value = -30138.285
$ws.Range("H96").Value = 23999
$ws.Range("J96").Value = 23999
$ws.Range("L96").Value = 23999
$ws.Range("N96").Value = -29491
$ws.Range("H97").Value = 23999
$ws.Range("J97").Value = 23999
$ws.Range("L97").Value = 23999
$ws.Range("N97").Value = -25981
$ws.Range("H101").Value = 26625
$ws.Range("J101").Value = 26625
$ws.Range("L101").Value = 26625
$ws.Range("N101").Value = -33115
$ws.Range("H104").Value = 25000
$ws.Range("J104").Value = 25000
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -30242
$ws.Range("H106").Value = 30440
$ws.Range("J106").Value = 30440
$ws.Range("L106").Value = 30440
$ws.Range("N106").Value = -32964
$ws.Range("H132").Value = 2735.2144
$ws.Range("I132").Value = 1316.25
$ws.Range("J132").Value = 4627.1665
$ws.Range("K132").Value = 3948.75
$ws.Range("L132").Value = 13881.4995
$ws.Range("M132").Value = -1418.75
$ws.Range("N132").Value = -18941.4995
$ws.Range("H134").Value = 4764.6
$ws.Range("I134").Value = 2460.4
$ws.Range("J134").Value = 5916.7
$ws.Range("K134").Value = 7381.200000000001
$ws.Range("L134").Value = 17750.1
$ws.Range("M134").Value = -4846.200000000001
$ws.Range("N134").Value = -22820.1
$ws.Range("H136").Value = 2082.5642
$ws.Range("I136").Value = 1153.3077
$ws.Range("J136").Value = 3941.077
$ws.Range("K136").Value = 3459.9231
$ws.Range("L136").Value = 11823.231
$ws.Range("M136").Value = -909.9231
$ws.Range("N136").Value = -16923.231

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1499.125
$ws.Range("J80").Value = 1541.8572
$ws.Range("L80").Value = 4625.571599999999
$ws.Range("N80").Value = -6497.571599999999
$ws.Range("H83").Value = 1499.125
$ws.Range("J83").Value = 1541.8572
$ws.Range("L83").Value = 13876.7148
$ws.Range("N83").Value = -23236.7148
$ws.Range("H122").Value = 2148.4
$ws.Range("J122").Value = 2523
$ws.Range("L122").Value = 22707
$ws.Range("N122").Value = -27607

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H132").Value = 3179.739
$ws.Range("I132").Value = 2954.5264
$ws.Range("K132").Value = 8863.5792
$ws.Range("M132").Value = -6333.5792

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 875.25
$ws.Range("I2").Value = 501
$ws.Range("K2").Value = 501
$ws.Range("M2").Value = -389

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 57143284
$ws.Range("I2").Value = 100000000
$ws.Range("K2").Value = 100000000
$ws.Range("M2").Value = -99999888
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H117").Value = 40900
$ws.Range("J117").Value = 40900
$ws.Range("L117").Value = 40900
$ws.Range("N117").Value = -50078
$ws.Range("H126").Value = 46338.273
$ws.Range("I126").Value = 67162.8
$ws.Range("K126").Value = 201488.4
$ws.Range("M126").Value = -199018.4
$ws.Range("H132").Value = 3406.077
$ws.Range("I132").Value = 3954.3845
$ws.Range("J132").Value = 2309.4614
$ws.Range("K132").Value = 11863.1535
$ws.Range("L132").Value = 6928.3842
$ws.Range("M132").Value = -9333.1535
$ws.Range("N132").Value = -11988.3842
$ws.Range("H136").Value = 2666.8333
$ws.Range("I136").Value = 1314.5714
$ws.Range("J136").Value = 4560
$ws.Range("K136").Value = 3943.7142
$ws.Range("L136").Value = 13680
$ws.Range("M136").Value = -1393.7142
$ws.Range("N136").Value = -18780
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
